$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5156438605213324
$ws.Range("D2").Value = 0.6094401512769929

$ws.Range("C3").Value = -0.510556004498761
$ws.Range("D3").Value = 0.6129610878900507

$ws.Range("C4").Value = 0.6803561824373288
$ws.Range("D4").Value = 0.5008859593471295

$ws.Range("C5").Value = -1.634238703352198
$ws.Range("D5").Value = 0.111435741756297

$ws.Range("C6").Value = -0.9148802137475198
$ws.Range("D6").Value = 0.3666976625278076

$ws.Range("C7").Value = 0.2163594693616468
$ws.Range("D7").Value = 0.8300004695558383

$ws.Range("C8").Value = -2.079964738118147
$ws.Range("D8").Value = 0.04513754033197004

$ws.Range("C9").Value = 0.7896262959206657
$ws.Range("D9").Value = 0.4352185543292824

$ws.Range("C10").Value = -1.518332588049684
$ws.Range("D10").Value = 0.138175341575141

$ws.Range("C11").Value = -1.806088412546325
$ws.Range("D11").Value = 0.07976400393105232
$ws.Range("G11").Value = "No"
